$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 267
$ws.Range("C3").Value = 165040
$ws.Range("C4").Value = 155991
$ws.Range("C7").Value = 5.48
$ws.Range("C8").Value = 64.92
